# Add a new "Modelo" column (F) to the results sheet, matching the
# formatting already used by the other header cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (E1) onto the new
# header cell (F1) so the new column matches the rest of the header row.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "Modelo"
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
